$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-6 down to 4-7
$ws.Rows("3:3").Insert()

# Row 2: Australian A-League Men | Melbourne City vs Macarthur FC
$ws.Range("A2").Value = "Australian A-League Men"
$ws.Range("B2").Value = "'2025-12-23"
$ws.Range("C2").Value = "05:15:00"
$ws.Range("D2").Value = "Melbourne City"
$ws.Range("E2").Value = "Macarthur FC"
$ws.Range("F2").Value = 1.72
$ws.Range("G2").Value = 1.74
$ws.Range("H2").Value = 5.5
$ws.Range("I2").Value = 5.7
$ws.Range("J2").Value = 4.1
$ws.Range("K2").Value = 4.2
$ws.Range("L2").Value = 1.35
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 4.3
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 2.14
$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 1.46
$ws.Range("S2").Value = 3.1
$ws.Range("T2").Value = 1.81
$ws.Range("U2").Value = 2.14
$ws.Range("V2").Value = 1.21
$ws.Range("W2").Value = 2.34
$ws.Range("X2").Value = 18
$ws.Range("Y2").Value = 21
$ws.Range("Z2").Value = 44
$ws.Range("AA2").Value = 140
$ws.Range("AB2").Value = 9.4
$ws.Range("AC2").Value = 9.2
$ws.Range("AD2").Value = 21
$ws.Range("AE2").Value = 70
$ws.Range("AF2").Value = 10.5
$ws.Range("AG2").Value = 9.6
$ws.Range("AH2").Value = 19.5
$ws.Range("AI2").Value = 70
$ws.Range("AJ2").Value = 17
$ws.Range("AK2").Value = 17
$ws.Range("AL2").Value = 32
$ws.Range("AM2").Value = 100
$ws.Range("AN2").Value = 9.8
$ws.Range("AO2").Value = 75

# Row 3: Friendly Matches | Tombense MG vs Desportiva
$ws.Range("A3").Value = "Friendly Matches"
$ws.Range("B3").Value = "'2025-12-23"
$ws.Range("C3").Value = "10:00:00"
$ws.Range("D3").Value = "Tombense MG"
$ws.Range("E3").Value = "Desportiva"
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 950
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.25
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.2
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 1.2
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4: Algerian Ligue 1 | MC Alger vs ES Ben Aknoun
$ws.Range("A4").Value = "Algerian Ligue 1"
$ws.Range("B4").Value = "'2025-12-23"
$ws.Range("C4").Value = "15:30:00"
$ws.Range("D4").Value = "MC Alger"
$ws.Range("E4").Value = "ES Ben Aknoun"
$ws.Range("F4").Value = 1.37
$ws.Range("G4").Value = 1.44
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = 16
$ws.Range("J4").Value = 4.3
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 1.46
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 2.9
$ws.Range("O4").Value = 1.43
$ws.Range("P4").Value = 1.64
$ws.Range("Q4").Value = 2.28
$ws.Range("R4").Value = 1.23
$ws.Range("S4").Value = 4.4
$ws.Range("T4").Value = 2.68
$ws.Range("U4").Value = 1.52
$ws.Range("V4").Value = 1.06
$ws.Range("W4").Value = 3.25
$ws.Range("X4").Value = 11.5
$ws.Range("Y4").Value = 32
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 5.8
$ws.Range("AC4").Value = 14
$ws.Range("AD4").Value = 65
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 6.8
$ws.Range("AG4").Value = 13.5
$ws.Range("AH4").Value = 48
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 13.5
$ws.Range("AK4").Value = 24
$ws.Range("AL4").Value = 380
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 10.5
$ws.Range("AO4").Value = 1000

# Row 5: Friendly Matches | Serra Branca EC vs Maguary
$ws.Range("A5").Value = "Friendly Matches"
$ws.Range("B5").Value = "'2025-12-23"
$ws.Range("C5").Value = "16:00:00"
$ws.Range("D5").Value = "Serra Branca EC"
$ws.Range("E5").Value = "Maguary"
$ws.Range("F5").Value = 1.05
$ws.Range("G5").Value = 600
$ws.Range("H5").Value = 1.05
$ws.Range("I5").Value = 870
$ws.Range("J5").Value = 1.05
$ws.Range("K5").Value = 32
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 1.3
$ws.Range("O5").Value = 1.01
$ws.Range("P5").Value = 1.3
$ws.Range("Q5").Value = 1.32
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 1.32
$ws.Range("T5").Value = 1.04
$ws.Range("U5").Value = 1.04
$ws.Range("V5").Value = 1.02
$ws.Range("W5").Value = 1.02
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# Row 6: Portuguese Primeira Liga | Guimaraes vs Sporting Lisbon
$ws.Range("A6").Value = "Portuguese Primeira Liga"
$ws.Range("B6").Value = "'2025-12-23"
$ws.Range("C6").Value = "17:45:00"
$ws.Range("D6").Value = "Guimaraes"
$ws.Range("E6").Value = "Sporting Lisbon"
$ws.Range("F6").Value = 8.4
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 1.46
$ws.Range("I6").Value = 1.47
$ws.Range("J6").Value = 4.9
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 1.38
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 3.85
$ws.Range("O6").Value = 1.31
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.37
$ws.Range("S6").Value = 3.4
$ws.Range("T6").Value = 2.18
$ws.Range("U6").Value = 1.76
$ws.Range("V6").Value = 3.1
$ws.Range("W6").Value = 1.12
$ws.Range("X6").Value = 15.5
$ws.Range("Y6").Value = 7.6
$ws.Range("Z6").Value = 8
$ws.Range("AA6").Value = 12
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 11
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 16
$ws.Range("AF6").Value = 80
$ws.Range("AG6").Value = 34
$ws.Range("AH6").Value = 29
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 340
$ws.Range("AK6").Value = 160
$ws.Range("AL6").Value = 140
$ws.Range("AM6").Value = 190
$ws.Range("AN6").Value = 250
$ws.Range("AO6").Value = 8

# Row 7: Friendly Matches | Necaxa vs Atletico San Luis
$ws.Range("A7").Value = "Friendly Matches"
$ws.Range("B7").Value = "'2025-12-23"
$ws.Range("C7").Value = "18:00:00"
$ws.Range("D7").Value = "Necaxa"
$ws.Range("E7").Value = "Atletico San Luis"
$ws.Range("F7").Value = 1.05
$ws.Range("G7").Value = 1000
$ws.Range("H7").Value = 1.05
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 1.04
$ws.Range("K7").Value = 32
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 1.32
$ws.Range("O7").Value = 1.02
$ws.Range("P7").Value = 1.32
$ws.Range("Q7").Value = 1.32
$ws.Range("R7").Value = 1.18
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 1.04
$ws.Range("U7").Value = 1.04
$ws.Range("V7").Value = 1.02
$ws.Range("W7").Value = 1.02
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000
